$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Cart_Page: add a new row describing the "drop down products" /
#    previewCartItem locator, widen column D, and move the selection.
# ---------------------------------------------------------------------
$cart = $wb.Worksheets.Item("Cart_Page")

# Enter the new row's values in the same order the original author must
# have used (C6 before A6) so new shared-strings come out in the same
# order as the target workbook (previewCartItem=78, drop down products=79).
$cart.Range("C6").Value = "previewCartItem"
$cart.Range("A6").Value = "drop down products"
$cart.Range("B6").Value = "class name"

# Column D grows from a bestFit 20.71 width to an explicit custom width
# of ~31.57; setting ColumnWidth explicitly also clears the old bestFit flag.
$cart.Columns.Item(4).ColumnWidth = 30.666666666666668

# Selection moves from B8 to A6, and this sheet is no longer the active tab
# (the new Menu sheet below becomes the active one).
$cart.Range("A6").Select()

# ---------------------------------------------------------------------
# 2) Add the new "Menu" sheet after Cart_Page (last tab).
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$menu = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$menu.Name = "Menu"

# Header row (all reuse already-existing shared strings).
$menu.Range("A1").Value = "name"
$menu.Range("B1").Value = "By"
$menu.Range("C1").Value = "locator"
$menu.Range("D1").Value = "expected"
$menu.Range("E1").Value = "Comments"

# Row 2 - top user cart / class name / navUser-item-cartLabel
$menu.Range("A2").Value = "top user cart"
$menu.Range("B2").Value = "class name"
$menu.Range("C2").Value = "navUser-item-cartLabel"

# Row 3 - top user cart / id / navbar-cart-icon
$menu.Range("A3").Value = "top user cart"
$menu.Range("B3").Value = "id"
$menu.Range("C3").Value = "navbar-cart-icon"

# Row 4 - top user cart / xpath / //*[@id="navbar-cart-icon"]
$menu.Range("A4").Value = "top user cart"
$menu.Range("B4").Value = "xpath"
$menu.Range("C4").Value = '//*[@id="navbar-cart-icon"]'

# Column widths matching the rest of the sheet layout.
$menu.Columns.Item(1).ColumnWidth = 20.833333333333332
$menu.Columns.Item(2).ColumnWidth = 16.333333333333332
$menu.Columns.Item(3).ColumnWidth = 80.16666666666667
$menu.Columns.Item(4).ColumnWidth = 30.666666666666668
$menu.Columns.Item(5).ColumnWidth = 22

# Match the rest of the workbook's page setup (portrait orientation).
$menu.PageSetup.Orientation = 1

# Carry over the styled-but-empty rows 9-23 (and D15:D20) from Cart_Page,
# matching formatting only (no values).
$cart.Range("A9:A23").Copy()
$menu.Range("A9").PasteSpecial(-4122)  # xlPasteFormats

$cart.Range("D15:D20").Copy()
$menu.Range("D15").PasteSpecial(-4122)  # xlPasteFormats

# Menu becomes the active sheet/tab, with A2 selected.
$menu.Range("A2").Select()
